# Sprints workbook — add sprint info for Gijs on the
# "sprint_6to9(vacation)" sheet (sprint 6..9 / vacation block).
#
# Rows touched (1-based, matching the sheet's own row numbers):
#   20 -> entry for User Story #28 gets a brand-new Gijs row (All / Gijs / hours / notes)
#   12 -> entry for User Story #21 switches from solo "Pauline" to "Pauline, Gijs"
#   10 -> entry for User Story #19 gets a Gijs contribution
#   13 -> entry for User Story #22 switches from solo "Pauline" to "Pauline, Gijs"
#   24 -> entry for User Story #31 gets a brand-new Gijs row (All / Gijs / hours)
#
# (cells are written in the same order the original author entered them so the
# shared-string table grows with matching indices)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sprint_6to9(vacation)")
$ws.Activate()

# --- Row 20 ------------------------------------------------------------
$ws.Range("C20").Value = "All"
$ws.Range("D20").Value = "Gijs"
$ws.Range("E20").Value = "6-8 hours"
$ws.Range("F20").Value = "14 hours"
$ws.Range("G20").Value = "y"
$ws.Range("H20").Value = "See sprint retrospectives."

# --- Row 12 ----------------------------------------------------------------
$ws.Range("D12").Value = "Pauline, Gijs"
$ws.Range("E12").Value = "P:1, G:2"
$ws.Range("F12").Value = "P:1 , G:3"
$ws.Range("G12").Value = "y"

# --- Row 10 --------------------------------------------------------------
$ws.Range("D10").Value = "Gijs"
$ws.Range("E10").Value = "2-3 hours"
$ws.Range("F10").Value = 4

# --- Row 13 ----------------------------------------------------------------
$ws.Range("C13").Value = "P: 1.2 G: rest"
$ws.Range("D13").Value = "Pauline, Gijs"
$ws.Range("E13").Value = "P:1, G:2-3"
$ws.Range("F13").Value = "P:1 , G:3.5"

# --- Row 24 ------------------------------------------------------------
$ws.Range("C24").Value = "All"
$ws.Range("D24").Value = "Gijs"
$ws.Range("E24").Value = "2 hours"
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = "y"

# Put the selection where the author appears to have left off editing.
$ws.Range("E25").Select()
